$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D: shared-string "N +- M" values updated
$ws.Range("D2").Value = "1064.8+-15.4"
$ws.Range("D3").Value = "933.1+-16.9"
$ws.Range("D4").Value = "207.4+-9.3"
$ws.Range("D5").Value = "1250.5+-15.0"
$ws.Range("D6").Value = "304.4+-14.1"
$ws.Range("D7").Value = "20.2+-4.5"
$ws.Range("D8").Value = "371.0+-12.6"
$ws.Range("D9").Value = "27.5+-3.3"
$ws.Range("D10").Value = "16.0+-4.4"
$ws.Range("D11").Value = "33.6+-5.8"
$ws.Range("D12").Value = "58.9+-6.4"
$ws.Range("D13").Value = "13.6+-2.9"

# Column E: Z_Score numeric values updated
$ws.Range("E2").Value = -13.32
$ws.Range("E3").Value = -5.85
$ws.Range("E4").Value = -6.31
$ws.Range("E5").Value = -11.62
$ws.Range("E6").Value = 10.18
$ws.Range("E7").Value = 6.89
$ws.Range("E8").Value = -9.44
$ws.Range("E9").Value = -4.76
$ws.Range("E10").Value = -2.49
$ws.Range("E11").Value = -2.01
$ws.Range("E12").Value = 9.609999999999999
$ws.Range("E13").Value = 2.91
$ws.Range("E14").Value = 3.21

# Column F: P-value numeric values updated
$ws.Range("F10").Value = 0.998
$ws.Range("F13").Value = 0.006
$ws.Range("F14").Value = 0.018
